{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\n// Find the index of the \"Ver no Jupiter...\" paragraph; the empty paragraph\n// immediately preceding it (right after the \"LOQ4205...\" requirement line)\n// is removed together with it and the \"\u00a9 2020 ...\" paragraph that follows.\nlet removeIdx = [];\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const text = paragraphs.items[i].text.trim();\n  if (text === \"Ver no Jupiter Salvar em pdf Salvar em docx\") {\n    removeIdx.push(i - 1); // preceding empty paragraph\n    removeIdx.push(i);     // this paragraph\n    removeIdx.push(i + 1); // the \"\u00a9 2020 ...\" paragraph right after\n  }\n}\n\n// Delete from the highest index down so earlier indices stay valid.\nremoveIdx = Array.from(new Set(removeIdx)).sort((a, b) => b - a);\nfor (const idx of removeIdx) {\n  if (idx >= 0 && idx < paragraphs.items.length) {\n    paragraphs.items[idx].delete();\n  }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the \"Ver no Jupiter...\" paragraph by its text. The paragraph that\n# immediately precedes it (a blank paragraph right after the \"LOQ4205...\"\n# requirement line) and the \"\u00a9 2020 ...\" paragraph right after it are\n# removed together with it.\n$jupiterIndex = 0\n$idx = 0\nforeach ($p in $d.Paragraphs) {\n    $idx = $idx + 1\n    $t = $p.Range.Text.Trim()\n    if ($t -eq \"Ver no Jupiter Salvar em pdf Salvar em docx\") {\n        $jupiterIndex = $idx\n    }\n}\n\nif ($jupiterIndex -gt 0) {\n    $startPara = $d.Paragraphs.Item($jupiterIndex - 1)\n    $endPara = $d.Paragraphs.Item($jupiterIndex + 1)\n    $r = $d.Range($startPara.Range.Start, $endPara.Range.End)\n    $r.Delete()\n}\n"}
